{"js": "// Fix wording in the \"Kurzfassung\" intro paragraph:\n//   \"dem tragen schwerer Taschen, Eink\u00e4ufen, und \u00e4hnlichem.\"\n// becomes\n//   \"dem Tragen schwerer Taschen, Eink\u00e4ufe, und \u00e4hnlichem.\"\n// (\"tragen\" -> capitalized \"Tragen\"; \"Eink\u00e4ufen\" -> \"Eink\u00e4ufe\")\n\nconst oldPhrase = \"dem tragen schwerer Taschen, Eink\u00e4ufen, und \u00e4hnlichem.\";\nconst newPhrase = \"dem Tragen schwerer Taschen, Eink\u00e4ufe, und \u00e4hnlichem.\";\n\nconst results = context.document.body.search(oldPhrase, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(newPhrase, Word.InsertLocation.replace);\n} else {\n  // Fallback: the phrase may have been split differently; patch the two\n  // smaller sub-edits independently so the script is resilient either way.\n  const r1 = context.document.body.search(\"dem tragen\", { matchCase: true });\n  r1.load(\"items\");\n  const r2 = context.document.body.search(\"Eink\u00e4ufen,\", { matchCase: true });\n  r2.load(\"items\");\n  await context.sync();\n\n  if (r1.items.length > 0) {\n    r1.items[0].insertText(\"dem Tragen\", Word.InsertLocation.replace);\n  }\n  if (r2.items.length > 0) {\n    r2.items[0].insertText(\"Eink\u00e4ufe,\", Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Fix wording in the \"Kurzfassung\" intro paragraph:\n#   \"dem tragen schwerer Taschen, Eink\u00e4ufen, und \u00e4hnlichem.\"\n# becomes\n#   \"dem Tragen schwerer Taschen, Eink\u00e4ufe, und \u00e4hnlichem.\"\n# (\"tragen\" -> capitalized \"Tragen\"; \"Eink\u00e4ufen\" -> \"Eink\u00e4ufe\")\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$found = $find.Execute(\n    \"dem tragen schwerer Taschen, Eink\u00e4ufen, und \u00e4hnlichem.\",  # FindText\n    $true,                                                     # MatchCase\n    $false,                                                    # MatchWholeWord\n    $false,                                                    # MatchWildcards\n    $false,                                                    # MatchSoundsLike\n    $false,                                                    # MatchAllWordForms\n    $true,                                                     # Forward\n    1,                                                         # Wrap (wdFindContinue)\n    $false,                                                    # Format\n    \"dem Tragen schwerer Taschen, Eink\u00e4ufe, und \u00e4hnlichem.\",   # ReplaceWith\n    2                                                          # Replace (wdReplaceAll)\n)\n\nif (-not $found) {\n    # Fallback: patch the two smaller sub-edits independently in case the\n    # longer phrase can't be matched as one run.\n    $f1 = $d.Content.Find\n    $f1.ClearFormatting()\n    $f1.Replacement.ClearFormatting()\n    $f1.Execute(\"dem tragen\", $true, $false, $false, $false, $false, $true, 1, $false, \"dem Tragen\", 2)\n\n    $f2 = $d.Content.Find\n    $f2.ClearFormatting()\n    $f2.Replacement.ClearFormatting()\n    $f2.Execute(\"Eink\u00e4ufen,\", $true, $false, $false, $false, $false, $true, 1, $false, \"Eink\u00e4ufe,\", 2)\n}\n"}
